# Scheduled-runner update: refresh market-price-derived Leve profit columns
# (currentAveragePrice* / LevePrice* / LeveProfit*) across multiple job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 3166.1
$ws.Range("I19").Value = 3517.8572
$ws.Range("K19").Value = 3517.8572
$ws.Range("M19").Value = -3342.8572
# Row 32
$ws.Range("H32").Value = 781.2727
$ws.Range("J32").Value = 524
$ws.Range("L32").Value = 524
$ws.Range("N32").Value = -1176
# Row 62
$ws.Range("H62").Value = 2044.5714
$ws.Range("I62").Value = 2102
$ws.Range("J62").Value = 1901
$ws.Range("K62").Value = 2102
$ws.Range("L62").Value = 1901
$ws.Range("M62").Value = -1478
$ws.Range("N62").Value = -3149
# Row 64
$ws.Range("H64").Value = 2908.75
$ws.Range("I64").Value = 2908.75
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2908.75
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2660.75
$ws.Range("N64").ClearContents()
# Row 65
$ws.Range("H65").Value = 2044.5714
$ws.Range("I65").Value = 2102
$ws.Range("J65").Value = 1901
$ws.Range("K65").Value = 10510
$ws.Range("L65").Value = 9505
$ws.Range("M65").Value = -7390
$ws.Range("N65").Value = -15745
# Row 67
$ws.Range("H67").Value = 2908.75
$ws.Range("I67").Value = 2908.75
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2908.75
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2050.75
$ws.Range("N67").ClearContents()
# Row 100
$ws.Range("H100").Value = 1524.375
$ws.Range("I100").Value = 1032.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1032.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -491.5
$ws.Range("N100").Value = -4082
# Row 132
$ws.Range("H132").Value = 1833038.1
$ws.Range("I132").Value = 1856778.9
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5570336.699999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5567806.699999999
$ws.Range("N132").Value = -20060
# Row 135
$ws.Range("H135").Value = 5214.207
$ws.Range("I135").Value = 638.95654
$ws.Range("J135").Value = 22752.666
$ws.Range("K135").Value = 5750.60886
$ws.Range("L135").Value = 204773.994
$ws.Range("M135").Value = -3215.60886
$ws.Range("N135").Value = -209843.994

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3310.0667
$ws.Range("I63").Value = 2651.9092
$ws.Range("K63").Value = 2651.9092
$ws.Range("M63").Value = -1965.9092
# Row 66
$ws.Range("H66").Value = 3310.0667
$ws.Range("I66").Value = 2651.9092
$ws.Range("K66").Value = 13250.454
$ws.Range("M66").Value = -9827.546
# Row 74
$ws.Range("H74").Value = 1192.5667
$ws.Range("I74").Value = 1055
$ws.Range("K74").Value = 1055
$ws.Range("M74").Value = -181
# Row 77
$ws.Range("H77").Value = 1192.5667
$ws.Range("I77").Value = 1055
$ws.Range("K77").Value = 5275
$ws.Range("M77").Value = -907
# Row 88
$ws.Range("H88").Value = 6659.4443
$ws.Range("I88").Value = 4721.3335
$ws.Range("K88").Value = 4721.3335
$ws.Range("M88").Value = -4315.3335
# Row 91
$ws.Range("H91").Value = 6659.4443
$ws.Range("I91").Value = 4721.3335
$ws.Range("K91").Value = 4721.3335
$ws.Range("M91").Value = -3317.3335
# Row 122
$ws.Range("H122").Value = 1366.1562
$ws.Range("I122").Value = 1360.5667
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 4081.7001
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1631.7001
$ws.Range("N122").Value = -9250
# Row 132
$ws.Range("H132").Value = 5399.397
$ws.Range("I132").Value = 6430.1665
$ws.Range("J132").Value = 3337.8572
$ws.Range("K132").Value = 19290.4995
$ws.Range("L132").Value = 10013.5716
$ws.Range("M132").Value = -16760.4995
$ws.Range("N132").Value = -15073.5716

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 418.0909
$ws.Range("I22").Value = 419.9
$ws.Range("K22").Value = 419.9
$ws.Range("M22").Value = -246.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 40
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -266
# Row 22
$ws.Range("H22").Value = 336.27274
$ws.Range("I22").Value = 344.33334
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 344.33334
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 5.666659999999979
$ws.Range("N22").Value = -1000
# Row 31
$ws.Range("H31").Value = 4364.3267
$ws.Range("I31").Value = 2269.2144
$ws.Range("J31").Value = 7157.8096
$ws.Range("K31").Value = 2269.2144
$ws.Range("L31").Value = 7157.8096
$ws.Range("M31").Value = -1974.2144
$ws.Range("N31").Value = -7747.8096
# Row 34
$ws.Range("H34").Value = 4364.3267
$ws.Range("I34").Value = 2269.2144
$ws.Range("J34").Value = 7157.8096
$ws.Range("K34").Value = 2269.2144
$ws.Range("L34").Value = 7157.8096
$ws.Range("M34").Value = -2067.2144
$ws.Range("N34").Value = -7561.8096
# Row 58
$ws.Range("H58").Value = 8773695
$ws.Range("I58").Value = 1717.037
$ws.Range("J58").Value = 30304914
$ws.Range("K58").Value = 1717.037
$ws.Range("L58").Value = 30304914
$ws.Range("M58").Value = -1514.037
$ws.Range("N58").Value = -30305320
# Row 62
$ws.Range("H62").Value = 30306556
$ws.Range("I62").Value = 3598
$ws.Range("J62").Value = 111114450
$ws.Range("K62").Value = 3598
$ws.Range("L62").Value = 111114450
$ws.Range("M62").Value = -2974
$ws.Range("N62").Value = -111115698
# Row 65
$ws.Range("H65").Value = 30306556
$ws.Range("I65").Value = 3598
$ws.Range("J65").Value = 111114450
$ws.Range("K65").Value = 17990
$ws.Range("L65").Value = 555572250
$ws.Range("M65").Value = -14870
$ws.Range("N65").Value = -555578490
# Row 132
$ws.Range("H132").Value = 4631563
$ws.Range("I132").Value = 1648.8334
$ws.Range("J132").Value = 8335494
$ws.Range("K132").Value = 4946.5002
$ws.Range("L132").Value = 25006482
$ws.Range("M132").Value = -2416.5002
$ws.Range("N132").Value = -25011542
# Row 136
$ws.Range("H136").Value = 8773695
$ws.Range("I136").Value = 1717.037
$ws.Range("J136").Value = 30304914
$ws.Range("K136").Value = 5151.111
$ws.Range("L136").Value = 90914742
$ws.Range("M136").Value = -2601.111
$ws.Range("N136").Value = -90919842

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 2293.95
$ws.Range("I86").Value = 1992
$ws.Range("J86").Value = 3199.8
$ws.Range("K86").Value = 5976
$ws.Range("L86").Value = 9599.400000000001
$ws.Range("M86").Value = -4790
$ws.Range("N86").Value = -11971.4
# Row 89
$ws.Range("H89").Value = 2293.95
$ws.Range("I89").Value = 1992
$ws.Range("J89").Value = 3199.8
$ws.Range("K89").Value = 17928
$ws.Range("L89").Value = 28798.2
$ws.Range("M89").Value = -12000
$ws.Range("N89").Value = -40654.2
# Row 107
$ws.Range("H107").Value = 378.82352
$ws.Range("I107").Value = 200.75
$ws.Range("K107").Value = 602.25
$ws.Range("M107").Value = 1317.75
# Row 122
$ws.Range("H122").Value = 1199.6666
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 7200
$ws.Range("M122").Value = -4750
# Row 131
$ws.Range("H131").Value = 3107041.2
$ws.Range("I131").Value = 100130
$ws.Range("J131").Value = 3175380.2
$ws.Range("K131").Value = 300390
$ws.Range("L131").Value = 9526140.600000001
$ws.Range("M131").Value = -295350
$ws.Range("N131").Value = -9536220.600000001
# Row 132
$ws.Range("H132").Value = 1716.3043
$ws.Range("I132").Value = 702.0833
$ws.Range("J132").Value = 2822.7273
$ws.Range("K132").Value = 6318.7497
$ws.Range("L132").Value = 25404.5457
$ws.Range("M132").Value = -3788.7497
$ws.Range("N132").Value = -30464.5457

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 22224214
$ws.Range("I122").Value = 33335010
$ws.Range("J122").Value = 2620
$ws.Range("K122").Value = 100005030
$ws.Range("L122").Value = 7860
$ws.Range("M122").Value = -100002580
$ws.Range("N122").Value = -12760
# Row 132
$ws.Range("H132").Value = 5246.8237
$ws.Range("I132").Value = 6499.1816
$ws.Range("J132").Value = 2950.8333
$ws.Range("K132").Value = 19497.5448
$ws.Range("L132").Value = 8852.499899999999
$ws.Range("M132").Value = -16967.5448
$ws.Range("N132").Value = -13912.4999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# Row 61
$ws.Range("I61").Value = 1160.8462
$ws.Range("J61").Value = 25642496
$ws.Range("K61").Value = 1160.8462
$ws.Range("L61").Value = 25642496
$ws.Range("M61").Value = -958.8462
$ws.Range("N61").Value = -25642900
# Row 113
$ws.Range("I113").Value = 1160.8462
$ws.Range("J113").Value = 25642496
$ws.Range("K113").Value = 1160.8462
$ws.Range("L113").Value = 25642496
$ws.Range("M113").Value = 1009.1538
$ws.Range("N113").Value = -25646836
# Row 132
$ws.Range("H132").Value = 7086.9165
$ws.Range("I132").Value = 9382.218000000001
$ws.Range("K132").Value = 28146.654
$ws.Range("M132").Value = -25616.654

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2593.9524
$ws.Range("I132").Value = 2415
$ws.Range("J132").Value = 2832.5557
$ws.Range("K132").Value = 7245
$ws.Range("L132").Value = 8497.667099999999
$ws.Range("M132").Value = -4715
$ws.Range("N132").Value = -13557.6671

